# Apply the "500k Heist" edit:
#  1. Remove the "Meta description: ..." paragraph that follows the title.
#  2. Insert a new bold paragraph ("Play 500k Heist for Free: Review and
#     Game Guide") right before the final "Prompt: ..." paragraph.
#  3. Replace the text of the final paragraph (still italic) with the
#     meta-description body copy (without the "Meta description:" label).

$d = $word.ActiveDocument

# --- Step 1: delete the "Meta description" paragraph (2nd paragraph) ---
$metaPara = $d.Paragraphs.Item(2)
$metaRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)
$metaRange.Delete()

# --- Step 2: insert a new bold paragraph right before the last paragraph ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 500k Heist for Free: Review and Game Guide</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($newParaXml)

# The inserted paragraph's trailing mark merged with the following (Prompt)
# paragraph, so split them apart again right after the inserted text.
$mergedPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$findRange = $d.Range($mergedPara.Range.Start, $mergedPara.Range.Start)
$findRange.Find.Execute("Play 500k Heist for Free: Review and Game Guide", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($findRange.End, $findRange.End)
$splitPoint.InsertParagraphAfter()

# --- Step 3: replace the final (Prompt) paragraph's text ---
$promptPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$promptRange = $d.Range($promptPara.Range.Start, $promptPara.Range.End)
$oldPrompt = 'Prompt: Design a cartoon-style feature image for the game "500K Heist" featuring a happy Maya warrior with glasses. Notes: The image should be eye-catching and energetic with vibrant colors. The warrior should have a confident look on his face and can be holding a bag of money to represent the 500K Heist. The glasses can add a modern touch to the ancient Maya warrior outfit. A background of a sparkling gemstone mine or casino can be added for context. Overall, the image should convey the excitement of high potential winnings and the theme of the game.'
$newPrompt = 'Read our review of 500k Heist and play for free. Enjoy high payouts, excellent graphics, and free launch modes, but with limited symbols and paylines.'
$promptRange.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newPrompt, 2)
